# Update TPM-derived NATMI ligand/receptor edge statistics (L1cam-Egfr) for
# sheet "Sheet1" with freshly recomputed values (re-run with new TPM data).
# Columns A-F (cluster/gene/cell-count/detection-rate) are unchanged;
# only the expression/specificity/edge-weight columns (G-T) are refreshed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.375839
$ws.Range("H2").Value = 16.127517
$ws.Range("I2").Value = 0.2354568587499626
$ws.Range("J2").Value = 0.2354568587499626
$ws.Range("M2").Value = 1.701929666666667
$ws.Range("N2").Value = 5.105789
$ws.Range("O2").Value = 0.02105622887134972
$ws.Range("P2").Value = 0.02105622887134972
$ws.Range("Q2").Value = 9.149299877323665
$ws.Range("R2").Value = 82.343698895913
$ws.Range("S2").Value = 0.004957833507168274
$ws.Range("T2").Value = 0.004957833507168274

$ws.Range("G3").Value = 5.375839
$ws.Range("H3").Value = 16.127517
$ws.Range("I3").Value = 0.2354568587499626
$ws.Range("J3").Value = 0.2354568587499626
$ws.Range("O3").Value = 0.7732971809418951
$ws.Range("P3").Value = 0.7732971809418953
$ws.Range("Q3").Value = 336.011155936533
$ws.Range("R3").Value = 3024.100403428797
$ws.Range("S3").Value = 0.1820781251047801
$ws.Range("T3").Value = 0.1820781251047801

$ws.Range("G4").Value = 5.375839
$ws.Range("H4").Value = 16.127517
$ws.Range("I4").Value = 0.2354568587499626
$ws.Range("J4").Value = 0.2354568587499626
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4338690000000001
$ws.Range("N4").Value = 1.301607
$ws.Range("O4").Value = 0.005367815805265532
$ws.Range("P4").Value = 0.005367815805265533
$ws.Range("Q4").Value = 2.332409891091
$ws.Range("R4").Value = 20.991689019819
$ws.Range("S4").Value = 0.001263889047856223
$ws.Range("T4").Value = 0.001263889047856223

$ws.Range("G5").Value = 5.375839
$ws.Range("H5").Value = 16.127517
$ws.Range("I5").Value = 0.2354568587499626
$ws.Range("J5").Value = 0.2354568587499626
$ws.Range("M5").Value = 15.972384
$ws.Range("N5").Value = 47.917152
$ws.Range("O5").Value = 0.1976099128607259
$ws.Range("P5").Value = 0.1976099128607259
$ws.Range("Q5").Value = 85.864964830176
$ws.Range("R5").Value = 772.7846834715841
$ws.Range("S5").Value = 0.04652860934004035
$ws.Range("T5").Value = 0.04652860934004036

$ws.Range("G6").Value = 5.375839
$ws.Range("H6").Value = 16.127517
$ws.Range("I6").Value = 0.2354568587499626
$ws.Range("J6").Value = 0.2354568587499626
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2157183333333333
$ws.Range("N6").Value = 0.647155
$ws.Range("O6").Value = 0.002668861520763652
$ws.Range("P6").Value = 0.002668861520763652
$ws.Range("Q6").Value = 1.159667029348333
$ws.Range("R6").Value = 10.437003264135
$ws.Range("S6").Value = 0.0006284017501176576
$ws.Range("T6").Value = 0.0006284017501176577

$ws.Range("I7").Value = 0.007131134316291014
$ws.Range("J7").Value = 0.007131134316291014
$ws.Range("M7").Value = 1.701929666666667
$ws.Range("N7").Value = 5.105789
$ws.Range("O7").Value = 0.02105622887134972
$ws.Range("P7").Value = 0.02105622887134972
$ws.Range("Q7").Value = 0.2770991113684444
$ws.Range("R7").Value = 2.493892002316
$ws.Range("S7").Value = 0.0001501547962761596
$ws.Range("T7").Value = 0.0001501547962761596

$ws.Range("I8").Value = 0.007131134316291014
$ws.Range("J8").Value = 0.007131134316291014
$ws.Range("O8").Value = 0.7732971809418951
$ws.Range("P8").Value = 0.7732971809418953
$ws.Range("S8").Value = 0.00551448606370585
$ws.Range("T8").Value = 0.005514486063705851

$ws.Range("I9").Value = 0.007131134316291014
$ws.Range("J9").Value = 0.007131134316291014
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.4338690000000001
$ws.Range("N9").Value = 1.301607
$ws.Range("O9").Value = 0.005367815805265532
$ws.Range("P9").Value = 0.005367815805265533
$ws.Range("Q9").Value = 0.070640236612
$ws.Range("R9").Value = 0.635762129508
$ws.Range("S9").Value = 0.00003827861549245832
$ws.Range("T9").Value = 0.00003827861549245833

$ws.Range("I10").Value = 0.007131134316291014
$ws.Range("J10").Value = 0.007131134316291014
$ws.Range("M10").Value = 15.972384
$ws.Range("N10").Value = 47.917152
$ws.Range("O10").Value = 0.1976099128607259
$ws.Range("P10").Value = 0.1976099128607259
$ws.Range("Q10").Value = 2.600538376832
$ws.Range("R10").Value = 23.404845391488
$ws.Range("S10").Value = 0.001409182830840399
$ws.Range("T10").Value = 0.0014091828308404

$ws.Range("I11").Value = 0.007131134316291014
$ws.Range("J11").Value = 0.007131134316291014
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2157183333333333
$ws.Range("N11").Value = 0.647155
$ws.Range("O11").Value = 0.002668861520763652
$ws.Range("P11").Value = 0.002668861520763652
$ws.Range("Q11").Value = 0.03512210853555556
$ws.Range("R11").Value = 0.31609897682
$ws.Range("S11").Value = 0.00001903200997614631
$ws.Range("T11").Value = 0.00001903200997614631

$ws.Range("G12").Value = 9.994147
$ws.Range("H12").Value = 29.982441
$ws.Range("I12").Value = 0.4377345486919088
$ws.Range("J12").Value = 0.4377345486919088
$ws.Range("M12").Value = 1.701929666666667
$ws.Range("N12").Value = 5.105789
$ws.Range("O12").Value = 0.02105622887134972
$ws.Range("P12").Value = 0.02105622887134972
$ws.Range("Q12").Value = 17.00933527232766
$ws.Range("R12").Value = 153.084017450949
$ws.Range("S12").Value = 0.009217038842153809
$ws.Range("T12").Value = 0.009217038842153809

$ws.Range("G13").Value = 9.994147
$ws.Range("H13").Value = 29.982441
$ws.Range("I13").Value = 0.4377345486919088
$ws.Range("J13").Value = 0.4377345486919088
$ws.Range("O13").Value = 0.7732971809418951
$ws.Range("P13").Value = 0.7732971809418953
$ws.Range("Q13").Value = 624.6736343982091
$ws.Range("R13").Value = 5622.062709583881
$ws.Range("S13").Value = 0.3384988925043258
$ws.Range("T13").Value = 0.3384988925043258

$ws.Range("G14").Value = 9.994147
$ws.Range("H14").Value = 29.982441
$ws.Range("I14").Value = 0.4377345486919088
$ws.Range("J14").Value = 0.4377345486919088
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.4338690000000001
$ws.Range("N14").Value = 1.301607
$ws.Range("O14").Value = 0.005367815805265532
$ws.Range("P14").Value = 0.005367815805265533
$ws.Range("Q14").Value = 4.336150564743001
$ws.Range("R14").Value = 39.025355082687
$ws.Range("S14").Value = 0.002349678428979202
$ws.Range("T14").Value = 0.002349678428979203

$ws.Range("G15").Value = 9.994147
$ws.Range("H15").Value = 29.982441
$ws.Range("I15").Value = 0.4377345486919088
$ws.Range("J15").Value = 0.4377345486919088
$ws.Range("M15").Value = 15.972384
$ws.Range("N15").Value = 47.917152
$ws.Range("O15").Value = 0.1976099128607259
$ws.Range("P15").Value = 0.1976099128607259
$ws.Range("Q15").Value = 159.630353636448
$ws.Range("R15").Value = 1436.673182728032
$ws.Range("S15").Value = 0.08650068602313726
$ws.Range("T15").Value = 0.08650068602313728

$ws.Range("G16").Value = 9.994147
$ws.Range("H16").Value = 29.982441
$ws.Range("I16").Value = 0.4377345486919088
$ws.Range("J16").Value = 0.4377345486919088
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.2157183333333333
$ws.Range("N16").Value = 0.647155
$ws.Range("O16").Value = 0.002668861520763652
$ws.Range("P16").Value = 0.002668861520763652
$ws.Range("Q16").Value = 2.155920733928333
$ws.Range("R16").Value = 19.403286605355
$ws.Range("S16").Value = 0.001168252893312679
$ws.Range("T16").Value = 0.001168252893312679

$ws.Range("G17").Value = 0.7761303333333333
$ws.Range("H17").Value = 2.328391
$ws.Range("I17").Value = 0.03399380269149206
$ws.Range("J17").Value = 0.03399380269149207
$ws.Range("M17").Value = 1.701929666666667
$ws.Range("N17").Value = 5.105789
$ws.Range("O17").Value = 0.02105622887134972
$ws.Range("P17").Value = 0.02105622887134972
$ws.Range("Q17").Value = 1.320919239499889
$ws.Range("R17").Value = 11.888273155499
$ws.Range("S17").Value = 0.0007157812896795609
$ws.Range("T17").Value = 0.000715781289679561

$ws.Range("G18").Value = 0.7761303333333333
$ws.Range("H18").Value = 2.328391
$ws.Range("I18").Value = 0.03399380269149206
$ws.Range("J18").Value = 0.03399380269149207
$ws.Range("O18").Value = 0.7732971809418951
$ws.Range("P18").Value = 0.7732971809418953
$ws.Range("Q18").Value = 48.511209219759
$ws.Range("R18").Value = 436.600882977831
$ws.Range("S18").Value = 0.02628731179082582
$ws.Range("T18").Value = 0.02628731179082583

$ws.Range("G19").Value = 0.7761303333333333
$ws.Range("H19").Value = 2.328391
$ws.Range("I19").Value = 0.03399380269149206
$ws.Range("J19").Value = 0.03399380269149207
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 0.4338690000000001
$ws.Range("N19").Value = 1.301607
$ws.Range("O19").Value = 0.005367815805265532
$ws.Range("P19").Value = 0.005367815805265533
$ws.Range("Q19").Value = 0.336738891593
$ws.Range("R19").Value = 3.030650024337
$ws.Range("S19").Value = 0.0001824724713684691
$ws.Range("T19").Value = 0.0001824724713684691

$ws.Range("G20").Value = 0.7761303333333333
$ws.Range("H20").Value = 2.328391
$ws.Range("I20").Value = 0.03399380269149206
$ws.Range("J20").Value = 0.03399380269149207
$ws.Range("M20").Value = 15.972384
$ws.Range("N20").Value = 47.917152
$ws.Range("O20").Value = 0.1976099128607259
$ws.Range("P20").Value = 0.1976099128607259
$ws.Range("Q20").Value = 12.396651718048
$ws.Range("R20").Value = 111.569865462432
$ws.Range("S20").Value = 0.006717512387670456
$ws.Range("T20").Value = 0.006717512387670458

$ws.Range("G21").Value = 0.7761303333333333
$ws.Range("H21").Value = 2.328391
$ws.Range("I21").Value = 0.03399380269149206
$ws.Range("J21").Value = 0.03399380269149207
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.2157183333333333
$ws.Range("N21").Value = 0.647155
$ws.Range("O21").Value = 0.002668861520763652
$ws.Range("P21").Value = 0.002668861520763652
$ws.Range("Q21").Value = 0.1674255419561111
$ws.Range("R21").Value = 1.506829877605
$ws.Range("S21").Value = 0.00009072475194775505
$ws.Range("T21").Value = 0.00009072475194775506

$ws.Range("G22").Value = 6.522593333333333
$ws.Range("H22").Value = 19.56778
$ws.Range("I22").Value = 0.2856836555503455
$ws.Range("J22").Value = 0.2856836555503455
$ws.Range("M22").Value = 1.701929666666667
$ws.Range("N22").Value = 5.105789
$ws.Range("O22").Value = 0.02105622887134972
$ws.Range("P22").Value = 0.02105622887134972
$ws.Range("Q22").Value = 11.10099509760222
$ws.Range("R22").Value = 99.90895587841999
$ws.Range("S22").Value = 0.006015420436071913
$ws.Range("T22").Value = 0.006015420436071913

$ws.Range("G23").Value = 6.522593333333333
$ws.Range("H23").Value = 19.56778
$ws.Range("I23").Value = 0.2856836555503455
$ws.Range("J23").Value = 0.2856836555503455
$ws.Range("O23").Value = 0.7732971809418951
$ws.Range("P23").Value = 0.7732971809418953
$ws.Range("Q23").Value = 407.68782800922
$ws.Range("R23").Value = 3669.19045208298
$ws.Range("S23").Value = 0.2209183654782576
$ws.Range("T23").Value = 0.2209183654782576

$ws.Range("G24").Value = 6.522593333333333
$ws.Range("H24").Value = 19.56778
$ws.Range("I24").Value = 0.2856836555503455
$ws.Range("J24").Value = 0.2856836555503455
$ws.Range("K24").Value = 3
$ws.Range("L24").Value = 1
$ws.Range("M24").Value = 0.4338690000000001
$ws.Range("N24").Value = 1.301607
$ws.Range("O24").Value = 0.005367815805265532
$ws.Range("P24").Value = 0.005367815805265533
$ws.Range("Q24").Value = 2.82995104694
$ws.Range("R24").Value = 25.46955942246
$ws.Range("S24").Value = 0.001533497241569179
$ws.Range("T24").Value = 0.001533497241569179

$ws.Range("G25").Value = 6.522593333333333
$ws.Range("H25").Value = 19.56778
$ws.Range("I25").Value = 0.2856836555503455
$ws.Range("J25").Value = 0.2856836555503455
$ws.Range("M25").Value = 15.972384
$ws.Range("N25").Value = 47.917152
$ws.Range("O25").Value = 0.1976099128607259
$ws.Range("P25").Value = 0.1976099128607259
$ws.Range("Q25").Value = 104.18136539584
$ws.Range("R25").Value = 937.63228856256
$ws.Range("S25").Value = 0.05645392227903741
$ws.Range("T25").Value = 0.05645392227903742

$ws.Range("G26").Value = 6.522593333333333
$ws.Range("H26").Value = 19.56778
$ws.Range("I26").Value = 0.2856836555503455
$ws.Range("J26").Value = 0.2856836555503455
$ws.Range("K26").Value = 3
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 0.2157183333333333
$ws.Range("N26").Value = 0.647155
$ws.Range("O26").Value = 0.002668861520763652
$ws.Range("P26").Value = 0.002668861520763652
$ws.Range("Q26").Value = 1.407042962877778
$ws.Range("R26").Value = 12.6633866659
$ws.Range("S26").Value = 0.0007624501154094147
$ws.Range("T26").Value = 0.0007624501154094147
